# Consider different semester and export timetable in excel file
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Modules offered in Semester 2 and Semester 3, appended below the
# existing Semester 1 timetable (rows 2-5).
$timetable = @(
    @{ Semester = 2; ModuleID = "SEG1201" },
    @{ Semester = 2; ModuleID = "OSS1014" },
    @{ Semester = 2; ModuleID = "WEB1201" },
    @{ Semester = 2; ModuleID = "PRG1203" },
    @{ Semester = 3; ModuleID = "SEG1201" },
    @{ Semester = 3; ModuleID = "OSS1014" },
    @{ Semester = 3; ModuleID = "WEB1201" },
    @{ Semester = 3; ModuleID = "NET1014" }
)

$row = 6
foreach ($entry in $timetable) {
    $ws.Cells.Item($row, 1).Value = $entry.Semester
    $ws.Cells.Item($row, 2).Value = $entry.ModuleID
    $row = $row + 1
}

# Select the newly added Semester 2 module column, as left by the export
$ws.Range("B6:B9").Select()
